# "Added black space around everything to avoid printing problems"
# Shrinks every shape on both slides (uniform scale + margin shift) so a
# safety border of blank space surrounds the printed content, and widens
# the crop on the big background photo on slide 1 to match.

$p = $ppt.ActivePresentation

# ---- Slide 1 ------------------------------------------------------------
$s1 = $p.Slides.Item(1)

$sh = $s1.Shapes.Item(1)   # Picture 3 (background photo)
$sh.PictureFormat.CropRight = 310.884525
$sh.Left   = 264.2018503937
$sh.Top    = 9.0931889764
$sh.Width  = 502.0072047244
$sh.Height = 475.5248425197

$sh = $s1.Shapes.Item(2)   # TextBox 4
$sh.Left   = 337.0392519685
$sh.Top    = 138.6405905512
$sh.Width  = 369.3901968504
$sh.Height = 114.4589370079

$sh = $s1.Shapes.Item(3)   # Rounded Rectangle 6
$sh.Left   = 283.9719291339
$sh.Top    = 41.9857874016
$sh.Width  = 468.2410629921
$sh.Height = 140.1513779528

$sh = $s1.Shapes.Item(4)   # TextBox 7
$sh.Left   = 518.6127165354
$sh.Top    = 91.8157874016
$sh.Width  = 216.1838976378
$sh.Height = 77.4364173228

$sh = $s1.Shapes.Item(5)   # TextBox 9
$sh.Left   = 300.3604330709
$sh.Top    = 59.4922440945
$sh.Width  = 201.8639763780
$sh.Height = 91.5157874016

$sh = $s1.Shapes.Item(6)   # Rectangle 10
$sh.Left   = 12.3922440945
$sh.Top    = 484.6179921260
$sh.Width  = 755.2553937008
$sh.Height = 46.8249212598

$sh = $s1.Shapes.Item(7)   # Picture 2
$sh.Left   = 634.2590944882
$sh.Top    = 486.1979133858
$sh.Width  = 131.9499606299
$sh.Height = 45.7425590551

$sh = $s1.Shapes.Item(8)   # Picture 4
$sh.Left   = 521.7343700787
$sh.Top    = 488.0917716535
$sh.Width  = 39.5364960630
$sh.Height = 39.8772047244

$sh = $s1.Shapes.Item(9)   # TextBox 11
$sh.Left   = 12.3922440945
$sh.Top    = 8.5736614173
$sh.Width  = 254.0515354331
$sh.Height = 373.1025590551

$sh = $s1.Shapes.Item(10)  # TextBox 12
$sh.Left   = 266.4436614173
$sh.Top    = 487.3866535433
$sh.Width  = 252.1690944882
$sh.Height = 42.2380708661

$sh = $s1.Shapes.Item(11)  # Picture 15
$sh.Left   = 12.3922440945
$sh.Top    = 491.9716141732
$sh.Width  = 124.8642913386
$sh.Height = 35.3064960630

$sh = $s1.Shapes.Item(12)  # TextBox 18
$sh.Left   = 519.7526377953
$sh.Top    = 54.1401181102
$sh.Width  = 210.7085433071
$sh.Height = 39.8914566929

# ---- Slide 2 ------------------------------------------------------------
$s2 = $p.Slides.Item(2)

$sh = $s2.Shapes.Item(1)   # Rectangle 4
$sh.Left   = 517.7352362205
$sh.Top    = 11.0894881890
$sh.Width  = 246.3245275591
$sh.Height = 171.8141338583

$sh = $s2.Shapes.Item(2)   # Picture 12
$sh.Left   = 303.5750000000
$sh.Top    = 358.0099606299
$sh.Width  = 453.1441338583
$sh.Height = 165.5441338583

$sh = $s2.Shapes.Item(3)   # Rectangle 13
$sh.Left   = 535.2354724409
$sh.Top    = 23.8250000000
$sh.Width  = 220.9628740157
$sh.Height = 146.3430314961

$sh = $s2.Shapes.Item(4)   # Rectangle 15
$sh.Left   = 520.4214566929
$sh.Top    = 271.8618503937
$sh.Width  = 239.3397244094
$sh.Height = 75.1265748031

$sh = $s2.Shapes.Item(5)   # TextBox 16
$sh.Left   = 530.0851574803
$sh.Top    = 199.3925590551
$sh.Width  = 210.6718503937
$sh.Height = 69.6872047244

$sh = $s2.Shapes.Item(6)   # TextBox 17
$sh.Left   = 37.3432677165
$sh.Top    = 39.1919291339
$sh.Width  = 379.2477559055
$sh.Height = 285.7173622047

$sh = $s2.Shapes.Item(7)   # Picture 1
$sh.Left   = 22.5302755906
$sh.Top    = 357.1131889764
$sh.Width  = 272.8675196850
$sh.Height = 165.8377559055

$sh = $s2.Shapes.Item(8)   # Rounded Rectangle 5
$sh.Left   = 520.8253937008
$sh.Top    = 192.6098031496
$sh.Width  = 238.0304330709
$sh.Height = 159.6363385827
